$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D8").Value = 22753300
$ws.Range("E8").Value = 20662500
$ws.Range("F8").Value = 16115400
$ws.Range("G8").Value = 13625600
$ws.Range("H8").Value = 13348800
$ws.Range("I8").Value = 12325300
$ws.Range("J8").Value = 10423600

$ws.Range("D9").Value = 29378500
$ws.Range("E9").Value = 12997400
$ws.Range("F9").Value = 9744600
$ws.Range("G9").Value = 7927500
$ws.Range("H9").Value = 7677500
$ws.Range("I9").Value = 7086100
$ws.Range("J9").Value = 6063900

$ws.Range("D10").Value = -6625200
$ws.Range("E10").Value = 7665100
$ws.Range("F10").Value = 6370800
$ws.Range("G10").Value = 5698000
$ws.Range("H10").Value = 5671300
$ws.Range("I10").Value = 5239200
$ws.Range("J10").Value = 4359800

$ws.Range("D14").Value = 106500
$ws.Range("F14").Value = 6900
$ws.Range("G14").Value = 5200

$ws.Range("D15").Value = 58300
$ws.Range("E15").Value = 43300
$ws.Range("F15").Value = 28100
$ws.Range("G15").Value = 21300
$ws.Range("H15").Value = 21000

$ws.Range("D17").Value = 20732400
$ws.Range("E17").Value = 18946700
$ws.Range("F17").Value = 14419900
$ws.Range("G17").Value = 12036100
$ws.Range("H17").Value = 11831900
$ws.Range("I17").Value = 10803000
$ws.Range("J17").Value = 9195400

$ws.Range("D18").Value = 2020900
$ws.Range("E18").Value = 1715800
$ws.Range("F18").Value = 1695400
$ws.Range("G18").Value = 1589500
$ws.Range("H18").Value = 1516900
$ws.Range("I18").Value = 1522300
$ws.Range("J18").Value = 1228200

$ws.Range("D20").Value = 1143600
$ws.Range("E20").Value = 53900
$ws.Range("F20").Value = -159900
$ws.Range("G20").Value = -160700
$ws.Range("H20").Value = -63900
$ws.Range("I20").Value = 31100
$ws.Range("J20").Value = 108900

$ws.Range("D21").Value = 4077200
$ws.Range("E21").Value = 2478100
$ws.Range("F21").Value = 2094900
$ws.Range("G21").Value = 1946200
$ws.Range("H21").Value = 1954000
$ws.Range("I21").Value = 1961100
$ws.Range("J21").Value = "NA"

$ws.Range("D22").Value = 904600
$ws.Range("E22").Value = 292800
$ws.Range("F22").Value = 234100
$ws.Range("G22").Value = 200700
$ws.Range("H22").Value = 155800
$ws.Range("I22").Value = 129600
$ws.Range("J22").Value = 119100

$ws.Range("D23").Value = 2259900
$ws.Range("E23").Value = 1476900
$ws.Range("F23").Value = 1301400
$ws.Range("G23").Value = 1228000
$ws.Range("H23").Value = 1297100
$ws.Range("I23").Value = 1423900
$ws.Range("J23").Value = 1218100

$ws.Range("D24").Value = 528200
$ws.Range("E24").Value = 408000
$ws.Range("F24").Value = 410200
$ws.Range("G24").Value = 323400
$ws.Range("H24").Value = 401100
$ws.Range("I24").Value = 411100
$ws.Range("J24").Value = 394000

$ws.Range("D26").Value = 1731600
$ws.Range("E26").Value = 1068900
$ws.Range("F26").Value = 891200
$ws.Range("G26").Value = 904600
$ws.Range("H26").Value = 896000
$ws.Range("I26").Value = 1012700
$ws.Range("J26").Value = 824100

$ws.Range("D27").Value = 2410400
$ws.Range("E27").Value = 1093400
$ws.Range("F27").Value = 914600
$ws.Range("G27").Value = 863800
$ws.Range("H27").Value = 823500
$ws.Range("I27").Value = 1071000
$ws.Range("J27").Value = 797500

$ws.Range("D29").Value = 192700
$ws.Range("E29").Value = "NA"
$ws.Range("F29").Value = "NA"
$ws.Range("G29").Value = "NA"
$ws.Range("H29").Value = "NA"
$ws.Range("I29").Value = "NA"
$ws.Range("J29").Value = "NA"

$ws.Range("D32").Value = -1143600
$ws.Range("E32").Value = -53900
$ws.Range("F32").Value = 159900
$ws.Range("G32").Value = 160700
$ws.Range("H32").Value = 63900
$ws.Range("I32").Value = -31100
$ws.Range("J32").Value = -108900

$ws.Range("D33").Value = 2603100
$ws.Range("E33").Value = 1093400
$ws.Range("F33").Value = 914600
$ws.Range("G33").Value = 863800
$ws.Range("H33").Value = 823500
$ws.Range("I33").Value = 1071000
$ws.Range("J33").Value = 797500

$ws.Range("D35").Value = 2603100
$ws.Range("E35").Value = 1093400
$ws.Range("F35").Value = 914600
$ws.Range("G35").Value = 863800
$ws.Range("H35").Value = 823500
$ws.Range("I35").Value = 1071000
$ws.Range("J35").Value = 797500

$ws.Range("D41").Value = 5013900
$ws.Range("E41").Value = 2256900
$ws.Range("F41").Value = 1520400
$ws.Range("G41").Value = 1835900
$ws.Range("H41").Value = 2826200
$ws.Range("I41").Value = 1888900
$ws.Range("J41").Value = 1336500

$ws.Range("D42").Value = 111700
$ws.Range("E42").Value = 6200
$ws.Range("F42").Value = 127000
$ws.Range("G42").Value = 7400
$ws.Range("H42").Value = 6500
$ws.Range("I42").Value = 88000
$ws.Range("J42").Value = 96100

$ws.Range("D43").Value = 2256000
$ws.Range("E43").Value = 1834100
$ws.Range("F43").Value = 2151200
$ws.Range("G43").Value = 1182900
$ws.Range("H43").Value = 1883800
$ws.Range("I43").Value = 935600
$ws.Range("J43").Value = 819300

$ws.Range("D44").Value = 1801900
$ws.Range("E44").Value = 1651500
$ws.Range("F44").Value = 1276500
$ws.Range("G44").Value = 890300
$ws.Range("H44").Value = 1891800
$ws.Range("I44").Value = 845400
$ws.Range("J44").Value = 1486600

$ws.Range("D45").Value = 187500
$ws.Range("E45").Value = 351700
$ws.Range("F45").Value = 441400
$ws.Range("G45").Value = 175100
$ws.Range("H45").Value = 985500
$ws.Range("I45").Value = 144800
$ws.Range("J45").Value = 107700

$ws.Range("D46").Value = 9371000
$ws.Range("E46").Value = 6100400
$ws.Range("F46").Value = 4485300
$ws.Range("G46").Value = 4091700
$ws.Range("H46").Value = 3805000
$ws.Range("I46").Value = 3902500
$ws.Range("J46").Value = 3102300

$ws.Range("D47").Value = 5618900
$ws.Range("E47").Value = 7524900
$ws.Range("F47").Value = 6327700
$ws.Range("G47").Value = 5686200
$ws.Range("H47").Value = 5229300
$ws.Range("I47").Value = 4539600
$ws.Range("J47").Value = 8044400

$ws.Range("D48").Value = 6036300
$ws.Range("E48").Value = 5287000
$ws.Range("F48").Value = 8305800
$ws.Range("G48").Value = 3911500
$ws.Range("H48").Value = 7649900
$ws.Range("I48").Value = 3188500
$ws.Range("J48").Value = 6016600

$ws.Range("D49").Value = 7969700
$ws.Range("E49").Value = 7927000
$ws.Range("F49").Value = 11206800
$ws.Range("G49").Value = 5251000
$ws.Range("H49").Value = 10684600
$ws.Range("I49").Value = 3511400
$ws.Range("J49").Value = 13061400

$ws.Range("D52").Value = 1443400
$ws.Range("E52").Value = 1380300
$ws.Range("F52").Value = 823200
$ws.Range("G52").Value = 515300
$ws.Range("H52").Value = 887100
$ws.Range("I52").Value = 164100
$ws.Range("J52").Value = 328300

$ws.Range("D54").Value = 30439300
$ws.Range("E54").Value = 28219600
$ws.Range("F54").Value = 21170700
$ws.Range("G54").Value = 19455700
$ws.Range("H54").Value = 18577400
$ws.Range("I54").Value = 15306100
$ws.Range("J54").Value = 13621100

$ws.Range("D57").Value = 3422000
$ws.Range("E57").Value = 3056100
$ws.Range("F57").Value = 2327900
$ws.Range("G57").Value = 1771200
$ws.Range("H57").Value = 1734800
$ws.Range("I57").Value = 1769100
$ws.Range("J57").Value = 1925900

$ws.Range("D58").Value = 702900
$ws.Range("E58").Value = 376600
$ws.Range("F58").Value = 610500
$ws.Range("G58").Value = 80300
$ws.Range("H58").Value = 197900
$ws.Range("I58").Value = 450100
$ws.Range("J58").Value = 288200

$ws.Range("D59").Value = 1306900
$ws.Range("E59").Value = 1030200
$ws.Range("F59").Value = 1498300
$ws.Range("G59").Value = 699300
$ws.Range("H59").Value = 2575400
$ws.Range("I59").Value = 290100
$ws.Range("J59").Value = 244100

$ws.Range("D60").Value = 5431700
$ws.Range("E60").Value = 4462900
$ws.Range("F60").Value = 3379700
$ws.Range("G60").Value = 2550800
$ws.Range("H60").Value = 2527500
$ws.Range("I60").Value = 2509200
$ws.Range("J60").Value = 2033900

$ws.Range("D61").Value = 6090400
$ws.Range("E61").Value = 6825300
$ws.Range("F61").Value = 4446300
$ws.Range("G61").Value = 4289400
$ws.Range("H61").Value = 3733400
$ws.Range("I61").Value = 1481300
$ws.Range("J61").Value = 1231900

$ws.Range("D62").Value = 1492100
$ws.Range("E62").Value = 2130700
$ws.Range("F62").Value = 1080200
$ws.Range("G62").Value = 713600
$ws.Range("H62").Value = 1574400
$ws.Range("I62").Value = 446100
$ws.Range("J62").Value = 416200

$ws.Range("D66").Value = 17494300
$ws.Range("E66").Value = 17259900
$ws.Range("F66").Value = 11782200
$ws.Range("G66").Value = 10638800
$ws.Range("H66").Value = 7067100
$ws.Range("I66").Value = 7276100
$ws.Range("J66").Value = 6161900

$ws.Range("D72").Value = 10440600
$ws.Range("E72").Value = 8730100
$ws.Range("F72").Value = 8127200
$ws.Range("G72").Value = 7609100
$ws.Range("H72").Value = 6767000
$ws.Range("I72").Value = 6646400
$ws.Range("J72").Value = 5921300

$ws.Range("D76").Value = 12945100
$ws.Range("E76").Value = 10959700
$ws.Range("F76").Value = 9388400
$ws.Range("G76").Value = 8816900
$ws.Range("H76").Value = 11510300
$ws.Range("I76").Value = 8030000
$ws.Range("J76").Value = 7459200

$ws.Range("D81").Value = 2603100
$ws.Range("E81").Value = 1093400
$ws.Range("F81").Value = 914600
$ws.Range("G81").Value = 863800
$ws.Range("H81").Value = 823500
$ws.Range("I81").Value = 1071000
$ws.Range("J81").Value = 797500

$ws.Range("D83").Value = 913600
$ws.Range("E83").Value = 709000
$ws.Range("F83").Value = 559900
$ws.Range("G83").Value = 517900
$ws.Range("H83").Value = 501500
$ws.Range("I83").Value = 408100
$ws.Range("J83").Value = "NA"

$ws.Range("D89").Value = 2075800
$ws.Range("E89").Value = 2592800
$ws.Range("F89").Value = 1900300
$ws.Range("G89").Value = 1932500
$ws.Range("H89").Value = 1487400
$ws.Range("I89").Value = 1592200
$ws.Range("J89").Value = 1098900

$ws.Range("D91").Value = -1077700
$ws.Range("E91").Value = -987000
$ws.Range("F91").Value = -904300
$ws.Range("G91").Value = -878500
$ws.Range("H91").Value = -847200
$ws.Range("I91").Value = -767700
$ws.Range("J91").Value = -1172000

$ws.Range("D94").Value = 1624900
$ws.Range("E94").Value = -1998700
$ws.Range("F94").Value = -1466700
$ws.Range("G94").Value = -807200
$ws.Range("H94").Value = -2856500
$ws.Range("I94").Value = -757300
$ws.Range("J94").Value = "NA"

$ws.Range("D96").Value = -643900
$ws.Range("E96").Value = -623000
$ws.Range("F96").Value = -553400
$ws.Range("G96").Value = -163000
$ws.Range("H96").Value = -853000
$ws.Range("I96").Value = -475100
$ws.Range("J96").Value = -342600

$ws.Range("D100").Value = -1114000
$ws.Range("E100").Value = 67100
$ws.Range("F100").Value = -710700
$ws.Range("G100").Value = -480400
$ws.Range("H100").Value = 1064600
$ws.Range("I100").Value = -176800
$ws.Range("J100").Value = "NA"

$ws.Range("D101").Value = 170400
$ws.Range("E101").Value = 75400
$ws.Range("F101").Value = -38400
$ws.Range("G101").Value = -218800
$ws.Range("H101").Value = -174500
$ws.Range("I101").Value = -105700
$ws.Range("J101").Value = "NA"

$ws.Range("D102").Value = 2757000
$ws.Range("E102").Value = 736500
$ws.Range("F102").Value = -315500
$ws.Range("G102").Value = 426100
$ws.Range("H102").Value = -479000
$ws.Range("I102").Value = 552400
$ws.Range("J102").Value = -44700
